$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ------------------
# The shared string "Ready for handoff" is used by every row of the
# generated localization-status report (Overview!E2:F3, zh-cn!C2:C3,
# de-de!C2:C3). Update every occurrence so the rendered text changes
# everywhere the status is shown.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width change: narrow the status columns ------------------------
# Overview columns E ("zh-cn") and F ("de-de"), and the "Status" column (C)
# on both the zh-cn and de-de sheets, shrink from ~17.22 characters wide to
# ~13.41 characters wide (report regenerated for archive with a narrower
# status column). ColumnWidth = 12.5 is the closest settable value that
# rounds to the target stored width in this engine's pixel grid.

$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe.Range("C1").ColumnWidth = 12.5
